$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for 12 new rows: insert blank rows 261-272 (old rows 261-269 shift to 273-281) ---
$ws.Rows("261:272").Insert()

# --- 2. Update existing row 260 quantity (15000 -> 0) ---
$ws.Cells.Item(260, 3).Value = 0

# --- 3. Fill the newly inserted rows (261-272) with the new order data ---
$data = @(
    @(261, "80266182", "12732-ROY-I", 5000),
    @(262, "80266182", "12732-ROY-I", 10000),
    @(263, "80266185", "10000-LDG-I", 1),
    @(264, "80266189", "10382-ARI-I", 1),
    @(265, "80266190", "19510-STM-I", 400),
    @(266, "80266191", "10119-ATE-I", 1),
    @(267, "80266191", "33642-ATE-I", 4),
    @(268, "80266191", "33380-ATE-I", 1),
    @(269, "80266191", "33503-ATE-I", 1),
    @(270, "80266194", "19876-WRN-I", 40000),
    @(271, "80266195", "22551-SHI-I", 500),
    @(272, "80266195", "60192-YAG-I", 6000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# --- 4. The previously-existing rows (now at 263-265, 275-280) already carry the
#        correct old values after the shift, except for three that changed:
#        row 263 -> 80266185 / 10000-LDG-I stays the same (already covered above)
#        row 264 -> 80266189 / 10382-ARI-I stays the same (already covered above)
#        row 265 -> 80266190 / 19510-STM-I stays the same (already covered above)
#        These were re-written above already as part of $data, so nothing else to do
#        for rows 261-272.

# --- 5. Append brand-new rows 273-274 and 281 (beyond the old data block) ---
$ws.Cells.Item(273, 1).Value = "80266198"
$ws.Cells.Item(273, 2).Value = "10493-ARI-I"
$ws.Cells.Item(273, 3).Value = 1

$ws.Cells.Item(274, 1).Value = "80266200"
$ws.Cells.Item(274, 2).Value = "21021-CTY-I"
$ws.Cells.Item(274, 3).Value = 1

# rows 275-280 already hold the shifted-down original data (84004808.. / 84004814..)
# and require no edits.

# --- 6. Add the final new row 281 ---
$ws.Rows("281:281").Insert()
$ws.Cells.Item(281, 1).Value = "84004819"
$ws.Cells.Item(281, 2).Value = "10255-ARI-I"
$ws.Cells.Item(281, 3).Value = 1

# --- 7. Refresh the view: selection now spans the full data table, and make sure
#        there is no stale "scrolled down" top-left cell recorded ---
$ws.Range("A1").Select()
$ws.Range("A2:C281").Select()
